# Update the "secante" (secant method) iteration table.
# The table now has 10 iteration rows (rows 2-11) instead of 4 (rows 2-5),
# and the previously-wrong values are replaced with the corrected,
# converged secant-method results.
# Numeric literals are passed as quoted strings cast to [double] so that
# scientific-notation values (e.g. 1.46e-06) parse correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("0", "1",                 "-75",                   "1.000005"),
    @("0", "7",                 "-72.29381928997699",     "1.000005"),
    @("1", "167.286012731271",  "-45.6396747455757",      "160.286012731271"),
    @("2", "441.742447644579",  "-17.9975833298854",      "274.456434913308"),
    @("3", "620.439221470155",  "-3.25561172458789",      "178.696773825576"),
    @("4", "659.902552942285",  "-0.202792517081662",     "39.4633314721297"),
    @("5", "662.524021092768",  "-0.0021885137978188",    "2.62146815048322"),
    @("6", "662.552620318726",  "-1.46057210770323e-06",  "0.0285992259574641"),
    @("7", "662.552639418047",  "-1.05160324892495e-11",  "1.9099321093563e-05"),
    @("8", "662.552639418184",  "0",                      "1.37561073643155e-10")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = [double]$row[0]
    $ws.Cells.Item($r, 2).Value = [double]$row[1]
    $ws.Cells.Item($r, 3).Value = [double]$row[2]
    $ws.Cells.Item($r, 4).Value = [double]$row[3]
    $r++
}

Write-Host "Updated secant table to $($r - 1) data rows"
